# Apply "feat: add 2022-Q3 data" to the workbook.
#
# Before:  总计 (sheet1), 2022-Q2 (sheet2), 2022-Q1 (sheet3)
# After:   总计 (sheet1), 2022-Q3 (new sheet2), 2022-Q2 (sheet3), 2022-Q1 (sheet4)
#
# The 总计 (summary) sheet gets a new row 2 for 2022-Q3, pushing the
# existing 2022-Q2 / 2022-Q1 rows down by one. A brand-new worksheet named
# "2022-Q3" is inserted right after 总计 and before 2022-Q2, holding the
# per-fund breakdown for that quarter.

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. 总计 sheet: insert the 2022-Q3 row, shifting the old rows down.
# ---------------------------------------------------------------------

# Row 4 did not exist before -- clone row 3's formatting (style s="2" on
# column A) onto it before writing new values.
$summary.Range("A3").Copy($summary.Range("A4"))

# Shift old row 3 (2022-Q1) down to row 4.
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(4, 2).Value = "2022-Q1"
$summary.Cells.Item(4, 3).Value = 20
$summary.Cells.Item(4, 4).Value = 2.37

# Shift old row 2 (2022-Q2) down to row 3.
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(3, 2).Value = "2022-Q2"
$summary.Cells.Item(3, 3).Value = 24
$summary.Cells.Item(3, 4).Value = 6.09

# New row 2: 2022-Q3 totals.
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 13
$summary.Cells.Item(2, 4).Value = 1.99

# ---------------------------------------------------------------------
# 2. Insert a new worksheet "2022-Q3" right after 总计 (before 2022-Q2).
# ---------------------------------------------------------------------

$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# Header row (B1:H1), cloning the bold/centered/bordered style used by the
# other quarter sheets' headers (style s="2" on 总计!B1).
$summary.Range("B1").Copy($q3.Range("B1:H1"))
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Column A (row index), cloning the style s="2" used on 总计!A2.
$summary.Range("A2").Copy($q3.Range("A2:A14"))

# Per-fund data rows. Columns B, D, E, F, G are stored as *text* in the
# target file (even though D/E/F/G look numeric), so each of those cells
# is pre-formatted as Text ("@") before the value is written, then reset
# to the default "Normal" style so no stray number-format style lingers.
# Column H (rank) is a genuine number and needs no special handling.

$rows = @(
    @(0,  "001480", "财通成长优选混合",             "20.31", "91.20", "5.14", "1.0439", 10),
    @(1,  "009693", "富国积极成长一年定期开放混合", "12.30", "97.80", "4.06", "0.4994", 6),
    @(2,  "011815", "恒越优势精选混合",             "3.22",  "90.91", "4.13", "0.1330", 1),
    @(3,  "002455", "民生加银鑫喜灵活配置混合",     "9.04",  "22.41", "0.86", "0.0777", 7),
    @(4,  "013028", "恒越品质生活混合",             "1.35",  "89.92", "4.05", "0.0547", 4),
    @(5,  "005265", "博时厚泽回报灵活配置混合A",    "1.70",  "72.26", "2.84", "0.0483", 9),
    @(6,  "012153", "博时研究慧选混合A",            "1.20",  "66.41", "3.15", "0.0378", 8),
    @(7,  "166109", "信澳量化先锋混合（LOF）A",     "0.79",  "88.99", "4.11", "0.0325", 4),
    @(8,  "005266", "博时厚泽回报灵活配置混合C",    "0.96",  "72.26", "2.84", "0.0273", 9),
    @(9,  "001250", "天弘新活力灵活配置混合",       "0.54",  "91.28", "3.95", "0.0213", 4),
    @(10, "012154", "博时研究慧选混合C",            "0.18",  "66.41", "3.15", "0.0057", 8),
    @(11, "000398", "华富灵活配置混合",             "0.12",  "94.04", "3.93", "0.0047", 4),
    @(12, "166110", "信澳量化先锋混合（LOF）C",     "0.11",  "88.99", "4.11", "0.0045", 4)
)

$r = 2
foreach ($d in $rows) {
    $q3.Cells.Item($r, 1).Value = $d[0]

    $cB = $q3.Cells.Item($r, 2)
    $cB.NumberFormat = "@"
    $cB.Value = $d[1]
    $cB.Style = "Normal"

    $q3.Cells.Item($r, 3).Value = $d[2]

    $cD = $q3.Cells.Item($r, 4)
    $cD.NumberFormat = "@"
    $cD.Value = $d[3]
    $cD.Style = "Normal"

    $cE = $q3.Cells.Item($r, 5)
    $cE.NumberFormat = "@"
    $cE.Value = $d[4]
    $cE.Style = "Normal"

    $cF = $q3.Cells.Item($r, 6)
    $cF.NumberFormat = "@"
    $cF.Value = $d[5]
    $cF.Style = "Normal"

    $cG = $q3.Cells.Item($r, 7)
    $cG.NumberFormat = "@"
    $cG.Value = $d[6]
    $cG.Style = "Normal"

    $q3.Cells.Item($r, 8).Value = $d[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3. Restore the original active sheet/tab (2022-Q1, the last sheet),
#    since adding the new sheet shifts Excel's notion of "active".
# ---------------------------------------------------------------------
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
